# Actualización automática hashcode
# Updates the hashcode values (column B) for the rows whose code (column A)
# matches the entries below, per the jue oct 3 02:26:54 CEST 2019 update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 34;  New = "82a122538dd440102d3a80b6a21db178" },
    @{ Row = 94;  New = "3976bbb9f1b4a382bc87fb541bc59088" },
    @{ Row = 95;  New = "5c73882c3c53b385a9b6cb3418168164" },
    @{ Row = 98;  New = "dfb77a4ff63d5cca57d0b52f6e8ac4ad" },
    @{ Row = 115; New = "662197525b2acd21c6124243032fe1bd" },
    @{ Row = 162; New = "496da3c040126f1aa643fcc0bd0ac7b4" },
    @{ Row = 180; New = "3628b7505f9fe43df36ba6974d4ef11f" },
    @{ Row = 213; New = "6d257cf1531177a2c618d10a50546c6c" },
    @{ Row = 227; New = "0bed25d524905a11baaf024e5fd8abc9" },
    @{ Row = 232; New = "72e804d3ceaaf08953cc162b25b3431f" },
    @{ Row = 420; New = "0841f66eec1f7caf51680bed6f5054c6" },
    @{ Row = 465; New = "89c67370eabfd551687d12306ce287f7" },
    @{ Row = 483; New = "0a8277e209a3872254017c3c03014112" },
    @{ Row = 513; New = "7ae2c5bb5dacbf5ba8bf260171240429" },
    @{ Row = 521; New = "68415814645a160ec90beea5ba8072f2" },
    @{ Row = 532; New = "7778078af76c44bcc102bff9c7d27ede" },
    @{ Row = 600; New = "eef16b95de2fdd043b7a987a50adf02f" },
    @{ Row = 626; New = "124054d4a6a4cbe2c5a28c761a12800c" },
    @{ Row = 674; New = "f8c310687a18f2145c2d5575eef9369e" },
    @{ Row = 708; New = "c162b077d372826d0847e23a22cd1573" },
    @{ Row = 737; New = "8ab5bc0100be605a7e27d1c9c2d71284" },
    @{ Row = 862; New = "cabac408ee7be64c2ee1efcd01eb2d8a" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.New
}
